$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$win = $wb.Windows.Item(1)
Write-Output "before: $($win.ScrollRow)"
$win.ScrollRow = 55
Write-Output "after: $($win.ScrollRow)"
